$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.983.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.742.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5148'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2759'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06197'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.741.46'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07237'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.17'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6513'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.633'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.82'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.003.67'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006801'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.964.97'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.296'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.683'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.386'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.97'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.516'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.27'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.786'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.80'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.953'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08246'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.663'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04687'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.657'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6257'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.733'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01616'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.921'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9999'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.35'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7638'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3873'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.017'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.351'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1134'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.67'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05230'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.79'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.577'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3437'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.56%  '
